$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Hoogte moet doorgestuurd worden door*DistanceSensor*") {
        $p.Range.Delete()
        break
    }
}
